# TrialsSetup update 2026-02-26 12:00
# The query table pulls "Days remaining" figures for each clinical trial.
# The scheduled refresh reduced the "Days remaining" value for the
# "REMASTER (CLOU)" trial (row 8) from 14 to 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 13
